# Add BOL in stage
$wb = $excel.ActiveWorkbook

$bol = $wb.Worksheets.Item("BOL")
$reroute = $wb.Worksheets.Item("Reroute Request")

# Update BOL test data: new order id / tracking number prefix for stage
$bol.Range("A2").Value = "58571277"
$bol.Range("A3").Value = ""
$bol.Range("A4").Value = ""
$bol.Range("C3").Value = "FCBTX"
$bol.Range("C4").Value = "FCBTX"

# Column A width tweak on BOL sheet
$bol.Columns.Item(1).ColumnWidth = 7.98828125

# Make BOL the active/selected sheet, with C11 selected
$bol.Activate()
$bol.Range("C11").Select()

# Reroute Request tab should no longer be the selected tab
$reroute.Range("I14").Select()
